$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows in column B (recalculated QoQ GVA series)
$updates = @{
    7  = "1.724740706211402"
    8  = "0.9852216748768466"
    23 = "1.62833067638352"
    24 = "0.9747899159663916"
    30 = "0.203404346429716"
    31 = "0.1816239316239334"
    40 = "-0.3816006600659966"
    41 = "0.6729475100942038"
    44 = "1.134157555941554"
    45 = "0.4041220448575383"
    46 = "0.2012477359629733"
    48 = "1.230123012301234"
    49 = "0.9583086346571812"
    52 = "-0.4417554979352672"
    53 = "0.5112375807851849"
    54 = "-0.5278310940499013"
    55 = "0.6753497346840218"
    56 = "0.287494010541458"
    57 = "-0.1051122790253219"
    59 = "-0.1524826074525977"
    60 = "-1.737138493843651"
    68 = "0.5792972459639122"
    69 = "-0.14162968558209"
    70 = "0.5673222390317647"
    72 = "-0.3945885005636995"
    73 = "-0.7074136955291455"
    74 = "0.1709888857224345"
    75 = "0.9009009009009035"
    76 = "-0.7612781954887239"
    77 = "-0.1420589070934801"
    78 = "-0.3793626707131937"
    80 = "0.5155131264916527"
    81 = "-0.09497578117580827"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = [double]$updates[$row]
}

# New row 82: appended date + value, A column uses the date style (style index 2)
$ws.Cells.Item(82, 1).Value = 45884
$ws.Cells.Item(82, 2).Value = 0.06654624964350926

# Apply the same style as the rest of column A (date format) to the new A82 cell
$ws.Range("A81").Copy()
$ws.Range("A82").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
